$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table from 45 to 58 rows by copying the last row's formatting
# (column A bold/bordered style, default style for B:F) down to the new rows.
$ws.Range("A45:F45").Copy($ws.Range("A46:F58"))

# Data rows 9-58: Day, Deaths, Detected Infections, Cumulative Infections,
# Current Asymptomatic Infections, Current Infectious Carriers
$data = @(
    @(7,0,0,30,9,21),
    @(8,0,0,37,16,21),
    @(9,0,0,43,22,21),
    @(10,0,0,47,26,21),
    @(11,0,21,48,27,0),
    @(12,0,21,53,24,8),
    @(13,0,21,63,27,15),
    @(14,0,21,72,30,21),
    @(15,0,21,82,35,26),
    @(16,0,21,88,40,27),
    @(17,0,29,98,46,23),
    @(18,0,36,107,46,25),
    @(19,0,42,122,51,29),
    @(20,0,47,145,67,31),
    @(21,1,48,177,92,37),
    @(22,1,52,205,111,42),
    @(23,1,61,233,130,42),
    @(24,1,71,261,145,45),
    @(25,1,78,298,161,59),
    @(26,1,85,345,182,78),
    @(27,2,94,410,217,99),
    @(28,2,103,488,269,116),
    @(29,3,116,546,297,133),
    @(30,3,137,627,352,138),
    @(31,3,163,708,385,160),
    @(32,3,193,799,421,185),
    @(33,3,219,921,469,233),
    @(34,6,249,1067,542,276),
    @(35,8,275,1219,626,318),
    @(36,9,323,1399,739,337),
    @(37,11,378,1572,820,374),
    @(38,12,452,1785,922,411),
    @(39,12,525,2007,998,484),
    @(40,13,593,2294,1139,562),
    @(41,15,660,2583,1270,653),
    @(42,17,752,2931,1455,724),
    @(43,18,863,3329,1667,799),
    @(44,21,1009,3781,1896,876),
    @(45,23,1155,4285,2132,998),
    @(46,27,1313,4827,2381,1133),
    @(47,34,1476,5441,2666,1299),
    @(48,38,1662,6166,3041,1463),
    @(49,43,1885,6941,3396,1660),
    @(50,48,2153,7801,3780,1868),
    @(51,53,2446,8806,4238,2122),
    @(52,62,2775,9941,4819,2347),
    @(53,71,3125,11100,5306,2669),
    @(54,90,3545,12327,5764,3018),
    @(55,106,4021,13691,6288,3382),
    @(56,118,4568,15165,6831,3766)
)

$startRow = 9
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
